$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# A1 keeps its existing (bold/centered/bordered) style but no longer holds a value.
$ws.Range("A1").Value = $null

$ws.Range("B1").Value = "Total"
$ws.Range("C1").Value = "Of which: male"
$ws.Range("D1").Value = "Of which: female"
$ws.Range("E1").Value = "By age and gender"
$ws.Range("F1").Value = "30 to 49"
$ws.Range("G1").Value = "Of which: male"
$ws.Range("H1").Value = "Of which: female"
$ws.Range("I1").Value = "50+"
$ws.Range("J1").Value = "Of which: male"
$ws.Range("K1").Value = "Of which: female"

# --- Row 2 (2019) ---
$ws.Range("A2").Value = 2019
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = 15
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = $null
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 17
$ws.Range("J2").Value = 13
$ws.Range("K2").Value = 4

# --- Row 3 (2018) ---
$ws.Range("A3").Value = 2018
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 15
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = $null
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 17
$ws.Range("J3").Value = 13
$ws.Range("K3").Value = 4

# --- Row 4 (2017) ---
$ws.Range("A4").Value = 2017
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = 15
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = $null
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 18
$ws.Range("J4").Value = 14
$ws.Range("K4").Value = 4
